$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.761.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.06%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.628.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.31%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.57%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'213.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.81%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.43%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.62%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -1.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.97%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.35%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.27%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.639.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.34%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'1.853.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.32%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -1.05%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₃0757"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.15%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'62.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.24%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'25.780.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.06%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.996"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.62%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.21%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'190.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.72%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.07%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.12%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.996"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.69%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -2.35%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'142.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.37%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.63%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.92%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.35%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.65%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0494"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.15%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.60%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.88%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.15%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.02%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.32%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.141.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.92%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.07%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -1.43%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.17%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.79%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'5.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.57%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'100.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.71%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.799"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.30%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.764.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.22%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0₆0109"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.20%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'55.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.48%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +2.07%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +5.17%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.44%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'7.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.46%  "
$ws.Range("E51").Style = "Normal"
